# Apply the changes described by the commit:
#  1. Slide 6's table switches to a different (built-in) table style.
#  2. The presentation's theme (theme1.xml, used by the slide master and
#     therefore by every slide) is recoloured from the custom "Integral"
#     palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Table style change on slide 6 (the table shape is Shapes.Item(2)).
# -----------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{0D360BFF-2A8C-4F9E-85F1-0939DBBFD58F}")
}

# -----------------------------------------------------------------
# 2) Recolour the theme to the default "Office" colour scheme.
# -----------------------------------------------------------------
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the twelve theme colour slots: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt $officeThemeColors[$i - 1]
}
